# Actualización 11 de Mayo - Mañana
# Updates grade-statistics numbers on the three partial-exam sheets.

$wb = $excel.ActiveWorkbook

# --- Sheet "1er Parcial" ---------------------------------------------------
$ws1 = $wb.Worksheets.Item("1er Parcial")

# Row 10: DESARROLLA APLICACIONES QUE SE EJECUTAN EN EL CLIENTE
$ws1.Cells.Item(10, 9).Value  = 6.2    # I10 Promedio
$ws1.Cells.Item(10, 10).Value = 2      # J10 Blancos
$ws1.Cells.Item(10, 11).Value = 7.14   # K10 Por_Blan

# Row 11: DESARROLLA APLICACIONES MÓVILES PARA IOS
$ws1.Cells.Item(11, 9).Value  = 6.2    # I11 Promedio
$ws1.Cells.Item(11, 10).Value = 0      # J11 Blancos
$ws1.Cells.Item(11, 11).Value = 0      # K11 Por_Blan

# --- Sheet "2o Parcial" -----------------------------------------------------
$ws2 = $wb.Worksheets.Item("2o Parcial")

# Row 9: CONSTRUYE PÁGINAS WEB
$ws2.Cells.Item(9, 5).Value  = 18             # E9 Aprobados
$ws2.Cells.Item(9, 6).Value  = 10             # F9 Reprobados
$ws2.Cells.Item(9, 7).Value  = 64.29000000000001  # G9 Por_Apro
$ws2.Cells.Item(9, 8).Value  = 35.71          # H9 Por_Repro
$ws2.Cells.Item(9, 9).Value  = 7.3            # I9 Promedio
$ws2.Cells.Item(9, 10).Value = 10             # J9 Blancos
$ws2.Cells.Item(9, 11).Value = 35.71          # K9 Por_Blan

# Row 10: DESARROLLA APLICACIONES QUE SE EJECUTAN EN EL CLIENTE
$ws2.Cells.Item(10, 5).Value  = 17   # E10 Aprobados
$ws2.Cells.Item(10, 6).Value  = 11   # F10 Reprobados
$ws2.Cells.Item(10, 7).Value  = 60.71  # G10 Por_Apro
$ws2.Cells.Item(10, 8).Value  = 39.29  # H10 Por_Repro
$ws2.Cells.Item(10, 9).Value  = 6.8    # I10 Promedio
$ws2.Cells.Item(10, 10).Value = 11     # J10 Blancos
$ws2.Cells.Item(10, 11).Value = 39.29  # K10 Por_Blan

# Row 11: DESARROLLA APLICACIONES MÓVILES PARA IOS
$ws2.Cells.Item(11, 5).Value  = 19              # E11 Aprobados
$ws2.Cells.Item(11, 6).Value  = 3               # F11 Reprobados
$ws2.Cells.Item(11, 7).Value  = 86.36           # G11 Por_Apro
$ws2.Cells.Item(11, 8).Value  = 13.64           # H11 Por_Repro
$ws2.Cells.Item(11, 9).Value  = 6.8             # I11 Promedio
$ws2.Cells.Item(11, 10).Value = 3               # J11 Blancos
$ws2.Cells.Item(11, 11).Value = 13.64           # K11 Por_Blan

# --- Sheet "3er Parcial" -----------------------------------------------------
$ws3 = $wb.Worksheets.Item("3er Parcial")

# Row 10: DESARROLLA APLICACIONES QUE SE EJECUTAN EN EL CLIENTE
$ws3.Cells.Item(10, 9).Value  = 6.4    # I10 Promedio
$ws3.Cells.Item(10, 10).Value = 2      # J10 Blancos
$ws3.Cells.Item(10, 11).Value = 7.14   # K10 Por_Blan

# Row 11: DESARROLLA APLICACIONES MÓVILES PARA IOS
$ws3.Cells.Item(11, 9).Value  = 6.5    # I11 Promedio
$ws3.Cells.Item(11, 10).Value = 0      # J11 Blancos
$ws3.Cells.Item(11, 11).Value = 0      # K11 Por_Blan
